$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks numeric need NumberFormat "@" (Text)
# applied first so Excel stores them as text rather than auto-converting
# them to numbers (matching the original inline-string / text cell type).
$textCells = @('D5','D7','D8','D9','D10','D11','D12','D14','D15','D16','D18','D20','D22','D23','D25','D26','D27','D28','D29','D30','D31','D32','D33','D34','D35','D36','D37','D39','D40','D41','D42','D43','D44','D45','D47','D48','D49','D50','D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '27.511.00'
$ws.Range('E2').Value = '  +2.14%  '
$ws.Range('D3').Value = '1.866.18'
$ws.Range('E3').Value = '  +2.68%  '
$ws.Range('D5').Value = '316.06'
$ws.Range('E5').Value = '  +2.34%  '
$ws.Range('E6').Value = '  -0.42%  '
$ws.Range('D7').Value = '0.4671'
$ws.Range('E7').Value = '  +0.86%  '
$ws.Range('D8').Value = '0.3737'
$ws.Range('E8').Value = '  +2.40%  '
$ws.Range('D9').Value = '0.07389'
$ws.Range('E9').Value = '  +2.33%  '
$ws.Range('D10').Value = '0.8899'
$ws.Range('E10').Value = '  +3.63%  '
$ws.Range('D11').Value = '0.07960'
$ws.Range('E11').Value = '  +5.36%  '
$ws.Range('D12').Value = '20.05'
$ws.Range('E12').Value = '  +1.61%  '
$ws.Range('D13').Value = '1.879.51'
$ws.Range('E13').Value = '  +5.88%  '
$ws.Range('D14').Value = '5.427'
$ws.Range('E14').Value = '  +2.04%  '
$ws.Range('D15').Value = '6.619'
$ws.Range('E15').Value = '  +2.31%  '
$ws.Range('D16').Value = '92.87'
$ws.Range('E16').Value = '  +1.19%  '
$ws.Range('E17').Value = '  -0.42%  '
$ws.Range('D18').Value = '0.000008950'
$ws.Range('E18').Value = '  +3.85%  '
$ws.Range('E19').Value = '  -0.29%  '
$ws.Range('D20').Value = '14.94'
$ws.Range('E20').Value = '  +3.50%  '
$ws.Range('D21').Value = '27.535.99'
$ws.Range('E21').Value = '  +3.29%  '
$ws.Range('D22').Value = '5.165'
$ws.Range('E22').Value = '  +0.59%  '
$ws.Range('D23').Value = '10.57'
$ws.Range('E23').Value = '  +0.67%  '
$ws.Range('D24').Value = '2.056.42'
$ws.Range('E24').Value = '  +5.40%  '
$ws.Range('D25').Value = '153.51'
$ws.Range('E25').Value = '  +1.05%  '
$ws.Range('D26').Value = '1.876'
$ws.Range('E26').Value = '  +1.48%  '
$ws.Range('D27').Value = '18.57'
$ws.Range('E27').Value = '  +2.43%  '
$ws.Range('D28').Value = '2.095'
$ws.Range('E28').Value = '  +1.22%  '
$ws.Range('D29').Value = '5.172'
$ws.Range('E29').Value = '  +1.76%  '
$ws.Range('D30').Value = '117.10'
$ws.Range('E30').Value = '  +1.57%  '
$ws.Range('D31').Value = '0.08918'
$ws.Range('E31').Value = '  +0.61%  '
$ws.Range('D32').Value = '0.7557'
$ws.Range('E32').Value = '  +5.76%  '
$ws.Range('D33').Value = '3.014'
$ws.Range('E33').Value = '  +1.47%  '
$ws.Range('D34').Value = '1.165'
$ws.Range('E34').Value = '  +3.38%  '
$ws.Range('D35').Value = '4.501'
$ws.Range('E35').Value = '  +2.14%  '
$ws.Range('D36').Value = '2.650'
$ws.Range('E36').Value = '  +10.49%  '
$ws.Range('D37').Value = '0.01976'
$ws.Range('E37').Value = '  +3.00%  '
$ws.Range('E38').Value = '  +0.90%  '
$ws.Range('D39').Value = '0.05291'
$ws.Range('E39').Value = '  +0.88%  '
$ws.Range('D40').Value = '2.990'
$ws.Range('E40').Value = '  +2.56%  '
$ws.Range('D41').Value = '7.193'
$ws.Range('E41').Value = '  +0.90%  '
$ws.Range('D42').Value = '0.5226'
$ws.Range('E42').Value = '  +1.62%  '
$ws.Range('D43').Value = '0.1651'
$ws.Range('E43').Value = '  +1.68%  '
$ws.Range('D44').Value = '8.378'
$ws.Range('E44').Value = '  +2.75%  '
$ws.Range('D45').Value = '0.4891'
$ws.Range('E45').Value = '  +1.87%  '
$ws.Range('E46').Value = '  +2.59%  '
$ws.Range('D47').Value = '1.004'
$ws.Range('E47').Value = '  -0.46%  '
$ws.Range('D48').Value = '103.86'
$ws.Range('E48').Value = '  +1.10%  '
$ws.Range('D49').Value = '1.666'
$ws.Range('E49').Value = '  +3.32%  '
$ws.Range('D50').Value = '0.06265'
$ws.Range('E50').Value = '  +0.04%  '
$ws.Range('D51').Value = '65.99'
$ws.Range('E51').Value = '  +2.95%  '
